$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.030683333333333
$ws.Range("H2").Value = 3.09205
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8686950000000001
$ws.Range("N2").Value = 2.606085
$ws.Range("O2").Value = 0.08920595562802248
$ws.Range("P2").Value = 0.08920595562802248
$ws.Range("Q2").Value = 0.8953494582500001
$ws.Range("R2").Value = 8.05814512425
$ws.Range("S2").Value = 0.08920595562802248
$ws.Range("T2").Value = 0.08920595562802248

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.030683333333333
$ws.Range("H3").Value = 3.09205
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.478549
$ws.Range("N3").Value = 13.435647
$ws.Range("O3").Value = 0.4599004752783479
$ws.Range("P3").Value = 0.4599004752783479
$ws.Range("Q3").Value = 4.615965811816667
$ws.Range("R3").Value = 41.54369230635
$ws.Range("S3").Value = 0.4599004752783479
$ws.Range("T3").Value = 0.4599004752783479

# Row 4 (Target cluster: sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.030683333333333
$ws.Range("H4").Value = 3.09205
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.390839
$ws.Range("N4").Value = 13.172517
$ws.Range("O4").Value = 0.4508935690936296
$ws.Range("P4").Value = 0.4508935690936296
$ws.Range("Q4").Value = 4.52556457665
$ws.Range("R4").Value = 40.73008118985
$ws.Range("S4").Value = 0.4508935690936296
$ws.Range("T4").Value = 0.4508935690936296
